# Update cryptocurrency price/volume snapshot values (D: Price, E: Volume(1h))
# for rows 2-51, as refreshed by the scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "67.101.21"
$ws.Range("E2").Value = "  +8.46%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.513.85"
$ws.Range("E3").Value = "  +11.71%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5: Solana
$ws.Range("D5").Value = "'191.53"
$ws.Range("E5").Value = "  +14.01%  "

# Row 6: BNB
$ws.Range("D6").Value = "'548.84"
$ws.Range("E6").Value = "  +6.28%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "3.509.13"
$ws.Range("E7").Value = "  +11.67%  "

# Row 8: XRP
$ws.Range("D8").Value = "'0.606"
$ws.Range("E8").Value = "  +3.76%  "

# Row 9: USDC
$ws.Range("E9").Value = "  -0.12%  "

# Row 10: Cardano
$ws.Range("D10").Value = "'0.631"
$ws.Range("E10").Value = "  +6.76%  "

# Row 11: Dogecoin
$ws.Range("E11").Value = "  +18.70%  "

# Row 12: Avalanche
$ws.Range("D12").Value = "'54.88"
$ws.Range("E12").Value = "  +6.38%  "

# Row 13: ShibaInu
$ws.Range("E13").Value = "  +9.86%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'9.33"
$ws.Range("E14").Value = "  +5.51%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.073.99"
$ws.Range("E15").Value = "  +11.63%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "3.513.35"
$ws.Range("E16").Value = "  +11.88%  "

# Row 17: TRON
$ws.Range("E17").Value = "  +5.23%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "67.041.92"
$ws.Range("E18").Value = "  +8.64%  "

# Row 19: Chainlink
$ws.Range("D19").Value = "'18.15"
$ws.Range("E19").Value = "  +7.58%  "

# Row 20: Uniswap
$ws.Range("D20").Value = "'11.89"
$ws.Range("E20").Value = "  +10.43%  "

# Row 21: Polygon
$ws.Range("D21").Value = "'0.993"
$ws.Range("E21").Value = "  +3.86%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "'422.44"
$ws.Range("E22").Value = "  +18.14%  "

# Row 23: PancakeSwap
$ws.Range("E23").Value = "  +6.41%  "

# Row 24: Litecoin
$ws.Range("E24").Value = "  +6.50%  "

# Row 25: Toncoin
$ws.Range("D25").Value = "'4.17"
$ws.Range("E25").Value = "  +7.67%  "

# Row 26: RenderToken
$ws.Range("D26").Value = "'11.15"
$ws.Range("E26").Value = "  +0.86%  "

# Row 27: ImmutableX
$ws.Range("D27").Value = "'2.89"
$ws.Range("E27").Value = "  +12.72%  "

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").Value = "'11.94"
$ws.Range("E28").Value = "  +8.11%  "

# Row 29: Filecoin
$ws.Range("D29").Value = "'8.88"
$ws.Range("E29").Value = "  +11.23%  "

# Row 30: EthereumClassic
$ws.Range("D30").Value = "'30.11"
$ws.Range("E30").Value = "  +8.74%  "

# Row 31: Bittensor
$ws.Range("D31").Value = "'652.41"
$ws.Range("E31").Value = "  +3.38%  "

# Row 32: NEARProtocol
$ws.Range("D32").Value = "'6.66"
$ws.Range("E32").Value = "  +5.79%  "

# Row 33: Cosmos
$ws.Range("D33").Value = "'11.67"
$ws.Range("E33").Value = "  +4.85%  "

# Row 34: Hedera
$ws.Range("E34").Value = "  +7.53%  "

# Row 35: OKB
$ws.Range("D35").Value = "'59.38"
$ws.Range("E35").Value = "  +6.20%  "

# Row 36: InjectiveProtocol
$ws.Range("D36").Value = "'38.49"
$ws.Range("E36").Value = "  +6.13%  "

# Row 37: PEPE
$ws.Range("D37").Value = "0.0₃0815"
$ws.Range("E37").Value = "  +19.61%  "

# Row 38: Dai
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.05%  "

# Row 39: TheGraph
$ws.Range("E39").Value = "  +5.83%  "

# Row 40: Kaspa
$ws.Range("E40").Value = "  +14.42%  "

# Row 41: Stacks
$ws.Range("D41").Value = "'3.30"
$ws.Range("E41").Value = "  +14.85%  "

# Row 42: FirstDigitalUSD
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.12%  "

# Row 43: Maker
$ws.Range("D43").Value = "2.995.00"
$ws.Range("E43").Value = "  +4.26%  "

# Row 44: Fetch.AI
$ws.Range("D44").Value = "'2.63"
$ws.Range("E44").Value = "  +5.72%  "

# Row 45: ThetaToken
$ws.Range("D45").Value = "'2.85"
$ws.Range("E45").Value = "  +13.98%  "

# Row 46: ApeXProtocol
$ws.Range("D46").Value = "'3.30"
$ws.Range("E46").Value = "  +12.82%  "

# Row 47: VeChain
$ws.Range("D47").Value = "'0.0416"
$ws.Range("E47").Value = "  +8.60%  "

# Row 48: WEMIXToken
$ws.Range("E48").Value = "  +4.10%  "

# Row 49: Stellar
$ws.Range("E49").Value = "  +7.05%  "

# Row 50: THORChain
$ws.Range("E50").Value = "  +16.66%  "

# Row 51: Monero
$ws.Range("D51").Value = "'140.09"
$ws.Range("E51").Value = "  +6.19%  "

